$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 222 (pushes existing rows 222-333 down to 223-334)
$ws.Rows.Item(222).EntireRow.Insert()

# Populate the newly inserted row 222 with the new price-report record
$ws.Range("A222").Value = 4
$ws.Range("B222").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C222").Value = "Los Lagos"
$ws.Range("D222").Value = 44572
$ws.Range("D222").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E222").Value = 10
$ws.Range("F222").Value = "Fruta"
$ws.Range("G222").Value = 100102
$ws.Range("H222").Value = "Cítricos"
$ws.Range("I222").Value = 100102005
$ws.Range("J222").Value = "Naranja"
$ws.Range("K222").Value = "Valencia"
$ws.Range("L222").Value = "Primera"
$ws.Range("M222").Value = 600
$ws.Range("N222").Value = 17000
$ws.Range("O222").Value = 18000
$ws.Range("P222").Value = 17500
$ws.Range("Q222").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R222").Value = "Región de O'Higgins"
$ws.Range("S222").Value = 1167
$ws.Range("T222").Value = 15
